# Roster sheet update: replace "Name"/"TL Name" columns with
# "Transport Status"/"Work Status" columns, and swap the per-employee
# Name values for Transport/Work status values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 3)
$ws.Range("B3").Value = "Transport Status"
$ws.Range("C3").Value = "Work Status"

# Employee rows (row 4 - row 10): fill Transport Status (col B) for every
# employee first, then Work Status (col C) for every employee, matching the
# order the values were authored in (keeps shared-string indices aligned).
$ws.Range("B4").Value = "Company Transport"
$ws.Range("B5").Value = "Self Transport"
$ws.Range("B6").Value = "Company Transport"
$ws.Range("B7").Value = "Company Transport"
$ws.Range("B8").Value = "Self Transport"
$ws.Range("B9").Value = "Company Transport"
$ws.Range("B10").Value = "Self Transport"

$ws.Range("C4").Value = "WFO"
$ws.Range("C5").Value = "WFH"
$ws.Range("C6").Value = "WFO"
$ws.Range("C7").Value = "WFH"
$ws.Range("C8").Value = "WFH"
$ws.Range("C9").Value = "WFO"
$ws.Range("C10").Value = "WFH"

# Move the active selection, matching the author's final cursor position.
$ws.Range("J12").Select() | Out-Null
